$d = $word.ActiveDocument

# Locate the paragraph that reads "You have the south integrated with the north"
$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -match "You have the south integrated with the north") {
        $targetIndex = $i
        break
    }
}

$target = $d.Paragraphs.Item($targetIndex)
$r = $target.Range

# Insert three new (empty) paragraphs immediately after the target paragraph.
$r.InsertParagraphAfter()
$r = $target.Range
$r.InsertParagraphAfter()
$r = $target.Range
$r.InsertParagraphAfter()

# Fill in the new paragraphs' text, in document order.
$d.Paragraphs.Item($targetIndex + 1).Range.Text = "will"
$d.Paragraphs.Item($targetIndex + 2).Range.Text = "will"
$d.Paragraphs.Item($targetIndex + 3).Range.Text = "wind wind wind wind wind wind wind wind wind wind wind wind wind wind wind wind wind wind wind "
